$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 (B5:AH5) values to 2 decimal places, matching the target dataset
$ws.Range("B5").Value = 5.08
$ws.Range("C5").Value = 4.16
$ws.Range("D5").Value = 0.34
$ws.Range("E5").Value = 11.15
$ws.Range("F5").Value = 9.369999999999999
$ws.Range("G5").Value = 4.18
$ws.Range("H5").Value = 19.74
$ws.Range("I5").Value = 6.35
$ws.Range("J5").Value = 2.9
$ws.Range("K5").Value = 4.65
$ws.Range("L5").Value = 5.28
$ws.Range("M5").Value = 4.73
$ws.Range("N5").Value = 1.17
$ws.Range("O5").Value = 3.91
$ws.Range("P5").Value = 5.97
$ws.Range("Q5").Value = 3.4
$ws.Range("R5").Value = 0.03
$ws.Range("S5").Value = 0.12
$ws.Range("T5").Value = 55.48
$ws.Range("U5").Value = 11.71
$ws.Range("V5").Value = 3.53
$ws.Range("W5").Value = 7.82
$ws.Range("X5").Value = 4.42
$ws.Range("Y5").Value = 0.5600000000000001
$ws.Range("Z5").Value = 9.609999999999999
$ws.Range("AA5").Value = 3.27
$ws.Range("AB5").Value = 3.47
$ws.Range("AC5").Value = 3.3
$ws.Range("AD5").Value = 5.2
$ws.Range("AE5").Value = 0.08
$ws.Range("AF5").Value = 18.14
$ws.Range("AG5").Value = 2.26
$ws.Range("AH5").Value = 4.64

# Remove row 6 entirely (trims the dataset, shrinking used range to A1:AH5)
$ws.Rows(6).Delete()
